$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected in the source file; unprotect to apply the data updates,
# then re-protect afterwards so the workbook keeps its protected state.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure note (shared string).
$ws.Range("A80").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-06 for illustrative purposes only and are subject to change."

# Refresh Weight (D) and Percent Change (E) figures for each holding row.
$ws.Range("D2").Value = 0.06368543776590888
$ws.Range("E2").Value = 0.01280249804839984
$ws.Range("D3").Value = 0.03871338828127497
$ws.Range("E3").Value = 0.01095537739945107
$ws.Range("D4").Value = 0.03209212451354953
$ws.Range("E4").Value = 0.01322676187771332
$ws.Range("D5").Value = 0.02984475132775059
$ws.Range("E5").Value = 0.00920553204484964
$ws.Range("D6").Value = 0.02739993694981467
$ws.Range("E6").Value = 0.009754748851937833
$ws.Range("D7").Value = 0.02486086057986931
$ws.Range("E7").Value = 0.02012442864398167
$ws.Range("D8").Value = 0.1864465901807476
$ws.Range("E8").Value = 0.01411535653443652
$ws.Range("D9").Value = 0.02439049757997625
$ws.Range("E9").Value = 0.004010295085892235
$ws.Range("D10").Value = 0.02211672787988082
$ws.Range("E10").Value = 0.01258804136070713
$ws.Range("D11").Value = 0.02215129192550672
$ws.Range("E11").Value = 0.006177325581395277
$ws.Range("D12").Value = 0.02005298605063773
$ws.Range("E12").Value = 0.001542614731970815
$ws.Range("D13").Value = 0.01992396272963695
$ws.Range("E13").Value = 0.01497946363855984
$ws.Range("D14").Value = 0.01674458346880453
$ws.Range("E14").Value = 0.02562096616467824
$ws.Range("D15").Value = 0.01632938089788978
$ws.Range("E15").Value = -0.004900255161215394
$ws.Range("D16").Value = 0.01504088378149788
$ws.Range("E16").Value = 0.008226652675760837
$ws.Range("D17").Value = 0.0144237814143871
$ws.Range("E17").Value = 0.003200568990042729
$ws.Range("D18").Value = 0.01419355172691323
$ws.Range("E18").Value = 0.0115713392174579
$ws.Range("D19").Value = 0.01367260527212017
$ws.Range("E19").Value = 0.01587200812646805
$ws.Range("D20").Value = 0.01347177080609754
$ws.Range("E20").Value = 0.009512875184517045
$ws.Range("D21").Value = 0.01251159266647858
$ws.Range("E21").Value = 0.01186387761473595
$ws.Range("D22").Value = 0.01324805146968334
$ws.Range("E22").Value = 0.002591121091725768
$ws.Range("D23").Value = 0.01165202904323636
$ws.Range("E23").Value = 0.01717505553448562
$ws.Range("D24").Value = 0.01302070211934496
$ws.Range("E24").Value = 0.005987878787878875
$ws.Range("D25").Value = 0.01179142946925953
$ws.Range("E25").Value = 0.01665411835566921
$ws.Range("D26").Value = 0.009407259997861291
$ws.Range("E26").Value = -0.004110393423370606
$ws.Range("D27").Value = 0.009643210902933057
$ws.Range("E27").Value = 0.01054009819967283
$ws.Range("D28").Value = 0.01039201964481466
$ws.Range("E28").Value = 0.005725611098877037
$ws.Range("D29").Value = 0.0099896799493269
$ws.Range("E29").Value = -0.003558718861209842
$ws.Range("D30").Value = 0.009765526589508721
$ws.Range("E30").Value = 0.01268686868686864
$ws.Range("D31").Value = 0.008748136273910813
$ws.Range("E31").Value = 0.005980650835532009
$ws.Range("D32").Value = 0.01023845427315256
$ws.Range("E32").Value = 0.01010459138450637
$ws.Range("D33").Value = 0.009494419788714724
$ws.Range("E33").Value = 0.01022320667916166
$ws.Range("D34").Value = 0.008948773456567992
$ws.Range("E34").Value = 0.01000000000000001
$ws.Range("D35").Value = 0.009298359580135859
$ws.Range("E35").Value = 0.005940762114911236
$ws.Range("D36").Value = 0.008416818589983117
$ws.Range("E36").Value = 0.002362669816893037
$ws.Range("D37").Value = 0.008610235201465023
$ws.Range("E37").Value = 0.006598845202089754
$ws.Range("D38").Value = 0.007941918072684826
$ws.Range("E38").Value = -0.01102930217307074
$ws.Range("D39").Value = 0.008831666050840063
$ws.Range("E39").Value = -0.01951463597698266
$ws.Range("D40").Value = 0.008168478289561427
$ws.Range("E40").Value = 0.0005506607929517404
$ws.Range("D41").Value = 0.007004822086822878
$ws.Range("E41").Value = 0.01070229592411498
$ws.Range("D42").Value = 0.007437977443993088
$ws.Range("E42").Value = 0.01002599331600429
$ws.Range("D43").Value = 0.008123892248970945
$ws.Range("E43").Value = 0.006945321378961866
$ws.Range("D44").Value = 0.007303193448721323
$ws.Range("E44").Value = 0.009681567203691177
$ws.Range("D45").Value = 0.007382659188322416
$ws.Range("E45").Value = 0.005066592556170813
$ws.Range("D46").Value = 0.007869238880855027
$ws.Range("E46").Value = 0.00890493381468116
$ws.Range("D47").Value = 0.007357880397622572
$ws.Range("E47").Value = 0.001029601029600924
$ws.Range("D48").Value = 0.007207392646461415
$ws.Range("E48").Value = 0.008562075044069362
$ws.Range("D49").Value = 0.006716314893196665
$ws.Range("E49").Value = 0.01110327811068013
$ws.Range("D50").Value = 0.007208260693269371
$ws.Range("E50").Value = 0.007947976878612817
$ws.Range("D51").Value = 0.006570167376075289
$ws.Range("E51").Value = -0.003405078190684252
$ws.Range("D52").Value = 0.006715604673081064
$ws.Range("E52").Value = 0.01890694586432606
$ws.Range("D53").Value = 0.005350324870858111
$ws.Range("E53").Value = -0.005973451327433543
$ws.Range("D54").Value = 0.006086112910620368
$ws.Range("E54").Value = 0.005601369223587893
$ws.Range("D55").Value = 0.006019707329811709
$ws.Range("E55").Value = -0.002064693737095613
$ws.Range("D56").Value = 0.005665875668219474
$ws.Range("E56").Value = 0.00760877286078987
$ws.Range("D57").Value = 0.006771396408828803
$ws.Range("E57").Value = 0.00121200820436318
$ws.Range("D58").Value = 0.005477217531512092
$ws.Range("E58").Value = -0.002766251728907432
$ws.Range("D59").Value = 0.005235032472092276
$ws.Range("E59").Value = 0.01451634784968125
$ws.Range("D60").Value = 0.004854117750091802
$ws.Range("E60").Value = 0.01261542463259202
$ws.Range("D61").Value = 0.004862798218171365
$ws.Range("E61").Value = 0.01865405212424109
$ws.Range("D62").Value = 0.004854749056861225
$ws.Range("E62").Value = 0.008127438231469331
$ws.Range("D63").Value = 0.004279312936532355
$ws.Range("E63").Value = 0.004130707383639409
$ws.Range("D64").Value = 0.00408329218462658
$ws.Range("E64").Value = -0.003401360544217913
$ws.Range("D65").Value = 0.003942037294968232
$ws.Range("E65").Value = -0.03587300316290987
$ws.Range("D66").Value = 0.003686989360121426
$ws.Range("E66").Value = 0.006677796327211993
$ws.Range("D67").Value = 0.003788787576690851
$ws.Range("E67").Value = 0.01430892276930762
$ws.Range("D68").Value = 0.003628396200584374
$ws.Range("E68").Value = 0.01445209278047832
$ws.Range("D69").Value = 0.00363072414429662
$ws.Range("E69").Value = 0.01201938751113918
$ws.Range("D70").Value = 0.003043372108694896
$ws.Range("E70").Value = -0.02338847689674839
$ws.Range("D71").Value = 0.002928474276659949
$ws.Range("E71").Value = 0.01867421180274853
$ws.Range("D72").Value = 0.002281937231424838
$ws.Range("E72").Value = 0.004461043676729837
$ws.Range("D73").Value = 0.001957366638595348
$ws.Range("E73").Value = 0.006994839542009412
$ws.Range("D74").Value = 0.001920316822564848
$ws.Range("E74").Value = -0.006328463703795029
$ws.Range("D75").Value = 0.001464947358445575
$ws.Range("E75").Value = -0.02149321266968318
$ws.Range("D76").Value = 0.001615435109606732
$ws.Range("E76").Value = 0.0115773533290997
$ws.Range("E77").Value = 0.009298521352495426

# Restore sheet protection (original password is unknown/unrecoverable from the
# legacy hash stored in the workbook, so we simply re-apply protection).
$ws.Protect($null)
